# Applies the LogicComponentClassDiagram content edits on slide 1:
#  - "AddressBook" -> "Erium" (the "Rectangle 62" box whose text reads
#    "AddressBookParser" across two paragraphs; only the first paragraph
#    changes).
#  - Inside the "Folded Corner 126" callout: resize/reposition the shape
#    and rename "AddCommand" -> "AddFriendCommand",
#    "FindCommand" -> "FindFriendCommand".

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "AddressBookParser" box: first paragraph "AddressBook" -> "Erium" ---
$addressBookShape = Get-ShapeById $s.Shapes 16
$tr = $addressBookShape.TextFrame.TextRange
$nameRange = $tr.Characters(1, 11)
$nameRange.Text = "Erium"

# --- "Folded Corner 126" callout: resize/reposition + text updates ---
$foldedCorner = Get-ShapeById $s.Shapes 127

$foldedCorner.Left = 493.1599
$foldedCorner.Top = 132.1279
$foldedCorner.Width = 103.6520
$foldedCorner.Height = 73.5157

$fcText = $foldedCorner.TextFrame.TextRange
$addCmd = $fcText.Characters(14, 10)
$addCmd.Text = "AddFriendCommand"

$fcText2 = $foldedCorner.TextFrame.TextRange
$findCmd = $fcText2.Characters(32, 11)
$findCmd.Text = "FindFriendCommand"
